# إضافة حدث جديد في Card11
# Fills the previously-blank placeholder cells on row 25 with the "nan"
# marker used throughout this sheet, then appends a new service-log row
# (row 26) describing the latest event for card 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card11")

# --- Row 25: the trailing columns (B:K and P) were left as empty cells;
# bring them in line with every other data row by stamping them "nan".
$ws.Range("B25:K25").Value = "nan"
$ws.Range("P25").Value = "nan"

# --- Row 26: brand new service event for card 11.
$cellA26 = $ws.Cells.Item(26, 1)
$cellA26.NumberFormat = "@"
$cellA26.Value = "11"
$cellA26.Style = "Normal"

$ws.Range("L26").Value = "29\11\2025"
$ws.Range("M26").Value = "1049.4 t"
$ws.Range("N26").Value = "تم سن الفلاتس"
$ws.Range("O26").Value = "الخبير"
